$wb = $excel.ActiveWorkbook

# --- Variables sheet: rename ZT3/ZR3/PF3 to lowercase ---
$ws1 = $wb.Worksheets.Item("Variables")
$ws1.Cells.Item(10,2).Value = "zt3"
$ws1.Cells.Item(11,2).Value = "zr3"
$ws1.Cells.Item(12,2).Value = "pf3"
$ws1.Range("B13").Select()

# --- Categories sheet: insert new category rows for disease variables ---
$ws2 = $wb.Worksheets.Item("Categories")

$ws2.Range("A55:A56").EntireRow.Insert()
$ws2.Cells.Item(55,1).Value = "casediab_fup5"
$ws2.Cells.Item(55,2).Value = 4
$ws2.Cells.Item(55,3).Value = 'inc. Diabetes (other types)'
$ws2.Cells.Item(56,1).Value = "casediab_fup5"
$ws2.Cells.Item(56,2).Value = 9
$ws2.Cells.Item(56,3).Value = 'incident (not verif.)'

$ws2.Rows.Item(54).Insert()
$ws2.Cells.Item(54,1).Value = "casediab_fup5"
$ws2.Cells.Item(54,2).Value = 1
$ws2.Cells.Item(54,3).Value = 'prevalent'

$ws2.Rows.Item(53).Insert()
$ws2.Cells.Item(53,1).Value = "casehf_fup5"
$ws2.Cells.Item(53,2).Value = 9
$ws2.Cells.Item(53,3).Value = 'incident (not verif.)'

$ws2.Rows.Item(52).Insert()
$ws2.Cells.Item(52,1).Value = "casehf_fup5"
$ws2.Cells.Item(52,2).Value = 1
$ws2.Cells.Item(52,3).Value = 'prevalent'

$ws2.Range("A51:A52").EntireRow.Insert()
$ws2.Cells.Item(51,1).Value = "casehyp_fup5"
$ws2.Cells.Item(51,2).Value = 3
$ws2.Cells.Item(51,3).Value = 'incident I15'
$ws2.Cells.Item(52,1).Value = "casehyp_fup5"
$ws2.Cells.Item(52,2).Value = 9
$ws2.Cells.Item(52,3).Value = 'incident (not verif.)'

$ws2.Rows.Item(50).Insert()
$ws2.Cells.Item(50,1).Value = "casehyp_fup5"
$ws2.Cells.Item(50,2).Value = 1
$ws2.Cells.Item(50,3).Value = 'prevalent'

$ws2.Rows.Item(45).Insert()
$ws2.Cells.Item(45,1).Value = "casestroke_fup5"
$ws2.Cells.Item(45,2).Value = 9
$ws2.Cells.Item(45,3).Value = 'incident (not verif.)'

$ws2.Rows.Item(43).Insert()
$ws2.Cells.Item(43,1).Value = "casestroke_fup5"
$ws2.Cells.Item(43,2).Value = 1
$ws2.Cells.Item(43,3).Value = 'prevalent'

$ws2.Rows.Item(42).Insert()
$ws2.Cells.Item(42,1).Value = "casemi_fup5"
$ws2.Cells.Item(42,2).Value = 9
$ws2.Cells.Item(42,3).Value = 'incident (not verif.)'

$ws2.Rows.Item(40).Insert()
$ws2.Cells.Item(40,1).Value = "casemi_fup5"
$ws2.Cells.Item(40,2).Value = 1
$ws2.Cells.Item(40,3).Value = 'prevalent'

# --- Cosmetic: column width + selection ---
$ws2.Columns.Item(1).ColumnWidth = 14.6
$ws2.Range("C66").Select()